$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.027.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.822.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.19%  '

$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4489'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.97%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3697'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07309'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8573'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.76'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.808.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.80%  '

$ws.Range("E13").Value = '  -1.14%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.332'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07100'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.95%  '

$ws.Range("E17").Value = '  -0.44%  '

$ws.Range("E18").Value = '  -1.04%  '

$ws.Range("E19").Value = '  -0.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.072.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.168'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.987'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.230'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.256'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08865'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.7554'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.183'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.79%  '

$ws.Range("E33").Value = '  +4.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.454'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.14%  '

$ws.Range("E35").Value = '  -0.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.091'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.17%  '

$ws.Range("E37").Value = '  -0.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05233'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5318'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.887'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.125'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1705'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5233'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.504'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.64'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.979'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.669'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06387'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9203'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
